$d = $word.ActiveDocument

function ReplaceText($findText, $replaceText) {
    $range = $d.Content
    $found = $range.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $found) {
        throw "Find failed for: $findText"
    }
}

# 1. Merge "Write Software Architecture Document section 1, 2, 3, 6: 5 hours." (bullet list item)
ReplaceText "Write Software Architecture Document section 1, 2, 3, 6: 5 hours." "Write Software Architecture Document section 1, 2, 3, 6: 5 hours."

# 2. Merge "carefully read the use case document..." (keep "All members are required to " separate)
ReplaceText "carefully read the use case document to understand all the use cases so that they can finish their tasks well." "carefully read the use case document to understand all the use cases so that they can finish their tasks well."

# 3. Merge "reading functionalities (ref: section 2.4.h)" (table cell, keep "UI design for " separate)
ReplaceText "reading functionalities (ref: section 2.4.h)" "reading functionalities (ref: section 2.4.h)"

# 4. Merge "social user profile (ref: section 2.4.i)" (table cell, keep "UI design for " separate)
ReplaceText "social user profile (ref: section 2.4.i)" "social user profile (ref: section 2.4.i)"

# 5. Merge "Refine the UI design for login, register, email verification " (keep trailing space)
ReplaceText "Refine the UI design for login, register, email verification " "Refine the UI design for login, register, email verification "

# 6. Merge "(ref: section 2.4.j)"
ReplaceText "(ref: section 2.4.j)" "(ref: section 2.4.j)"

# 7. Merge "Write Software Architecture Document section 1, 2, 3, 6" (table cell, row 10, no trailing hours)
ReplaceText "Write Software Architecture Document section 1, 2, 3, 6" "Write Software Architecture Document section 1, 2, 3, 6"

# 8. Merge "Prepare for Software Architecture Document section 4" (table cell, row 11)
ReplaceText "Prepare for Software Architecture Document section 4" "Prepare for Software Architecture Document section 4"

# 9. Add a new table row (row 12) after row 11, by duplicating row 11 and editing its cell contents.
$tables = $d.Tables
$lastTable = $tables.Item($tables.Count)
$lastRow = $lastTable.Rows.Item($lastTable.Rows.Count)
$newRow = $lastTable.Rows.Add()

# Copy cell contents: Cell1 = "12", Cell2 = task text, Cell3 = date, Cell4 = assignee
$newRow.Cells.Item(1).Range.Text = "12"
$newRow.Cells.Item(2).Range.Text = "Jira backlogs and assign tasks for assignees on Jira"
$newRow.Cells.Item(3).Range.Text = "13/11/2024"
$newRow.Cells.Item(4).Range.Text = "Pham Thanh Vinh"
